$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3
$ws.Range("D2").Value = "CALLE 10A NO. 1 - 21 SUR, FACA "
$ws.Range("E2").Value = 3174325821
$ws.Range("F2").Value = "carnesfaca.10@gmail.com"

$ws.Range("D3").Value = "VEREDA CHIGUALA, VILLAPINZÓN "
$ws.Range("E3").Value = 3112764625
$ws.Range("F3").Value = "jhon.cartertorres@gmail.com"
